$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reproduce the extra (unused) cell-format slot the source workbook carries
# in styles.xml (cellXfs index 1) so the style table lines up with the
# authored file before we add the boolean-formatted cell below.
$ws.Range("A1:B3").Locked = $false
$ws.Range("A1:B3").Locked = $true

# Row 4: new "new layer" boolean flag field
$ws.Range("A4").Value = "שכבה חדשה"
$ws.Range("B4").NumberFormat = '"TRUE";"TRUE";"FALSE"'
$ws.Range("B4").Value = $false

# Row 5: new "new layer location" text field
$ws.Range("A5").Value = "מיקום שכבה חדשה"
$ws.Range("B5").Value = "C:\Users\dpere\Documents\JTMT\Projects\תחזיות_דמוגרפיות\קבצי עבודה\142_מתחם_אנגל\בהת\For_approval\Reference_tabels\shp"

# Row 6: trailing formatted-but-empty row, keeps the used range in sync
$ws.Range("A6:B6").NumberFormat = "General"

# Match the authored row heights
$ws.Rows.Item(4).RowHeight = 13.8
$ws.Rows.Item(5).RowHeight = 14.5
$ws.Rows.Item(6).RowHeight = 13.8

# Move the active selection onto the newly added block
[void]$ws.Range("A4:B5").Select()
